$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.320.53"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").Value = "1.841.11"
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +1.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.26"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4731"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3698"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07452"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8855"
$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.51"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").Value = "1.840.09"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07346"
$ws.Range("E13").Value = "  +3.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.482"
$ws.Range("E14").Value = "  +2.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.34"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.591"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008858"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.84"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").Value = "27.324.36"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.319"

$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").Value = "2.071.13"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.62"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.188"
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.292"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.74"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08947"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7652"
$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.572"
$ws.Range("E34").Value = "  +1.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.934"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("E36").Value = "  +0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05349"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01967"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.000"
$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.370"
$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.418"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5376"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.577"
$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4983"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.63"
$ws.Range("E47").Value = "  +1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("E48").Value = "  +1.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.683"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.30"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06325"
$ws.Range("E51").Value = "  +0.43%  "

